$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting of the existing header cells (bold, bordered, centered/top-aligned)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for column I (I0) and column J (IF), rows 2-14
$dataI = @(6, 8, 9, 10, 8, 6, 5, 5, 6, 8, 8, 2, 4)
$dataJ = @(8, 9, 9, 10, 8, 7, 6, 6, 7, 8, 8, 2, 4)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
